$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.96712006850333
$ws.Range("C2").Value = 5.549828386289488
$ws.Range("E2").Value = 10.01842814961051
$ws.Range("F2").Value = 53.46964407262195
$ws.Range("G2").Value = 3.786699805609232
$ws.Range("I2").Value = 40.53025755385092
$ws.Range("J2").Value = 11.31122530482046
$ws.Range("K2").Value = 14.62324852237586
$ws.Range("L2").Value = 10.97786872625303
$ws.Range("M2").Value = 16.56554641144104
$ws.Range("B3").Value = 15.92528125975034
$ws.Range("C3").Value = 5.368998870191458
$ws.Range("E3").Value = 10.04676511012211
$ws.Range("F3").Value = 53.26135235750908
$ws.Range("G3").Value = 3.789746167686348
$ws.Range("I3").Value = 40.39695498833275
$ws.Range("J3").Value = 11.30807998222485
$ws.Range("K3").Value = 14.5857403505834
$ws.Range("L3").Value = 11.00779728507985
$ws.Range("M3").Value = 16.59969193537903
$ws.Range("B4").Value = 15.90439230393854
$ws.Range("C4").Value = 5.256257724735225
$ws.Range("E4").Value = 10.06530145201455
$ws.Range("F4").Value = 53.14042297342213
$ws.Range("G4").Value = 3.791714337022929
$ws.Range("I4").Value = 40.31864716642426
$ws.Range("J4").Value = 11.3061161191312
$ws.Range("K4").Value = 14.56657524734508
$ws.Range("L4").Value = 11.02786579998475
$ws.Range("M4").Value = 16.62386230730682
$ws.Range("B5").Value = 15.89709524003171
$ws.Range("C5").Value = 5.209960054319952
$ws.Range("E5").Value = 10.07314201078152
$ws.Range("F5").Value = 53.09291833361137
$ws.Range("G5").Value = 3.792541033523725
$ws.Range("I5").Value = 40.28763738465576
$ws.Range("J5").Value = 11.30530750189093
$ws.Range("K5").Value = 14.55974368079249
$ws.Range("L5").Value = 11.03646993430003
$ws.Range("M5").Value = 16.63451793202637
$ws.Range("B6").Value = 15.89595718685115
$ws.Range("C6").Value = 5.202253249723467
$ws.Range("E6").Value = 10.07446127983714
$ws.Range("F6").Value = 53.08513810407518
$ws.Range("G6").Value = 3.792679797277828
$ws.Range("I6").Value = 40.28254297420582
$ws.Range("J6").Value = 11.30517272500589
$ws.Range("K6").Value = 14.55866856352305
$ws.Range("L6").Value = 11.03792439176726
$ws.Range("M6").Value = 16.63633597171687
$ws.Range("B7").Value = 15.90428896238902
$ws.Range("C7").Value = 5.255634670670475
$ws.Range("E7").Value = 10.0654060299352
$ws.Range("F7").Value = 53.13977509160813
$ws.Range("G7").Value = 3.791725386213913
$ws.Range("I7").Value = 40.31822529458005
$ws.Range("J7").Value = 11.30610524772539
$ws.Range("K7").Value = 14.56647914534434
$ws.Range("L7").Value = 11.02798011259921
$ws.Range("M7").Value = 16.62400274927644
$ws.Range("B8").Value = 15.95170316574197
$ws.Range("C8").Value = 5.487877112029066
$ws.Range("E8").Value = 10.02796319459781
$ws.Range("F8").Value = 53.39639513857782
$ws.Range("G8").Value = 3.787729967646257
$ws.Range("I8").Value = 40.48356684087562
$ws.Range("J8").Value = 11.31014749383835
$ws.Range("K8").Value = 14.60951770736148
$ws.Range("L8").Value = 10.98783717991471
$ws.Range("M8").Value = 16.57665471359536
$ws.Range("B9").Value = 16.08233942707229
$ws.Range("C9").Value = 5.926694219902755
$ws.Range("E9").Value = 9.963523837749797
$ws.Range("F9").Value = 53.95372623192129
$ws.Range("G9").Value = 3.780666154154484
$ws.Range("I9").Value = 40.83538082150802
$ws.Range("J9").Value = 11.31782117688103
$ws.Range("K9").Value = 14.72425186525396
$ws.Range("L9").Value = 10.92252120791221
$ws.Range("M9").Value = 16.50922641852971
$ws.Range("B10").Value = 16.20060421814579
$ws.Range("C10").Value = 6.235316035821937
$ws.Range("E10").Value = 9.921605850134812
$ws.Range("F10").Value = 54.39453536239336
$ws.Range("G10").Value = 3.775940991190958
$ws.Range("I10").Value = 41.1099445932315
$ws.Range("J10").Value = 11.32331273252736
$ws.Range("K10").Value = 14.82652803653588
$ws.Range("L10").Value = 10.88267355811526
$ws.Range("M10").Value = 16.47516771680518
$ws.Range("B11").Value = 16.25907423267379
$ws.Range("C11").Value = 6.372068661094415
$ws.Range("E11").Value = 9.903703340956618
$ws.Range("F11").Value = 54.60148947140926
$ws.Range("G11").Value = 3.773891104081758
$ws.Range("I11").Value = 41.2381725854973
$ws.Range("J11").Value = 11.32578091043337
$ws.Range("K11").Value = 14.87683040990137
$ws.Range("L11").Value = 10.86630646323346
$ws.Range("M11").Value = 16.46302886626013
$ws.Range("B12").Value = 16.28187061955351
$ws.Range("C12").Value = 6.423282613389726
$ws.Range("E12").Value = 9.897090969119379
$ws.Range("F12").Value = 54.68074647296852
$ws.Range("G12").Value = 3.773129099660902
$ws.Range("I12").Value = 41.28719170154411
$ws.Range("J12").Value = 11.32671136263605
$ws.Range("K12").Value = 14.8964091971123
$ws.Range("L12").Value = 10.86036115665057
$ws.Range("M12").Value = 16.45891380313224
$ws.Range("B13").Value = 16.27693215220109
$ws.Range("C13").Value = 6.41227888529327
$ws.Range("E13").Value = 9.898507650196549
$ws.Range("F13").Value = 54.66363816098892
$ws.Range("G13").Value = 3.773292578753894
$ws.Range("I13").Value = 41.27661428306519
$ws.Range("J13").Value = 11.32651115779384
$ws.Range("K13").Value = 14.89216917781454
$ws.Range("L13").Value = 10.86163036036685
$ws.Range("M13").Value = 16.45977864716638
$ws.Range("B14").Value = 16.26093667589063
$ws.Range("C14").Value = 6.376293768153879
$ws.Range("E14").Value = 9.903155996246381
$ws.Range("F14").Value = 54.60799240972949
$ws.Range("G14").Value = 3.77382812852348
$ws.Range("I14").Value = 41.24219625663591
$ws.Range("J14").Value = 11.32585754268525
$ws.Range("K14").Value = 14.87843061484399
$ws.Range("L14").Value = 10.86581228057692
$ws.Range("M14").Value = 16.46268066919371
$ws.Range("B15").Value = 16.25122378085993
$ws.Range("C15").Value = 6.354176118279152
$ws.Range("E15").Value = 9.906024957381506
$ws.Range("F15").Value = 54.57402231677072
$ws.Range("G15").Value = 3.774158020730642
$ws.Range("I15").Value = 41.22117387908479
$ws.Range("J15").Value = 11.32545664184686
$ws.Range("K15").Value = 14.87008402686828
$ws.Range("L15").Value = 10.86840670276133
$ws.Range("M15").Value = 16.46452094287529
$ws.Range("B16").Value = 16.19687566178715
$ws.Range("C16").Value = 6.22630146293523
$ws.Range("E16").Value = 9.922799221373115
$ws.Range("F16").Value = 54.38113669647295
$ws.Range("G16").Value = 3.776076953678503
$ws.Range("I16").Value = 41.1016301539723
$ws.Range("J16").Value = 11.32315083831871
$ws.Range("K16").Value = 14.82331558639775
$ws.Range("L16").Value = 10.88377855550028
$ws.Range("M16").Value = 16.47602846356744
$ws.Range("B17").Value = 16.16471972082364
$ws.Range("C17").Value = 6.14688630479047
$ws.Range("E17").Value = 9.933387823476943
$ws.Range("F17").Value = 54.26442855848048
$ws.Range("G17").Value = 3.777279612488261
$ws.Range("I17").Value = 41.02913468130874
$ws.Range("J17").Value = 11.32172871184679
$ws.Range("K17").Value = 14.79558316433264
$ws.Range("L17").Value = 10.89365906481605
$ws.Range("M17").Value = 16.48394668055682
$ws.Range("B18").Value = 16.14666517953215
$ws.Range("C18").Value = 6.100868604841817
$ws.Range("E18").Value = 9.939587921002078
$ws.Range("F18").Value = 54.19790841559551
$ws.Range("G18").Value = 3.777980730904956
$ws.Range("I18").Value = 40.98775144547556
$ws.Range("J18").Value = 11.32090791577588
$ws.Range("K18").Value = 14.77998848153717
$ws.Range("L18").Value = 10.89950774401801
$ws.Range("M18").Value = 16.4888168363345
$ws.Range("B19").Value = 16.14062839968159
$ws.Range("C19").Value = 6.085230913978473
$ws.Range("E19").Value = 9.941706051906868
$ws.Range("F19").Value = 54.17549122364937
$ws.Range("G19").Value = 3.778219731200336
$ws.Range("I19").Value = 40.97379425632747
$ws.Range("J19").Value = 11.32062951933393
$ws.Range("K19").Value = 14.77476995159553
$ws.Range("L19").Value = 10.9015164755602
$ws.Range("M19").Value = 16.49052004734655
$ws.Range("B20").Value = 16.16809728081453
$ws.Range("C20").Value = 6.155375756037818
$ws.Range("E20").Value = 9.932249289062788
$ws.Range("F20").Value = 54.27678974541744
$ws.Range("G20").Value = 3.777150617110561
$ws.Range("I20").Value = 41.03681954289148
$ws.Range("J20").Value = 11.32188039134487
$ws.Range("K20").Value = 14.79849854417657
$ws.Range("L20").Value = 10.89259012495524
$ws.Range("M20").Value = 16.4830710929278
$ws.Range("B21").Value = 16.26561729481132
$ws.Range("C21").Value = 6.386879328553231
$ws.Range("E21").Value = 9.901786140332923
$ws.Range("F21").Value = 54.62431311333658
$ws.Range("G21").Value = 3.773670438656684
$ws.Range("I21").Value = 41.25229326325258
$ws.Range("J21").Value = 11.32604963818481
$ws.Range("K21").Value = 14.88245167511617
$ws.Range("L21").Value = 10.86457709896816
$ws.Range("M21").Value = 16.46181520959244
$ws.Range("B22").Value = 16.33316196386678
$ws.Range("C22").Value = 6.534831832189682
$ws.Range("E22").Value = 9.882849293287951
$ws.Range("F22").Value = 54.85659763928683
$ws.Range("G22").Value = 3.771478924725908
$ws.Range("I22").Value = 41.39580227389887
$ws.Range("J22").Value = 11.32875009251551
$ws.Range("K22").Value = 14.94040481292037
$ws.Range("L22").Value = 10.84774082101375
$ws.Range("M22").Value = 16.45073027070872
$ws.Range("B23").Value = 16.29676929133194
$ws.Range("C23").Value = 6.456187605482631
$ws.Range("E23").Value = 9.892867508378215
$ws.Range("F23").Value = 54.73216344550816
$ws.Range("G23").Value = 3.772641010537643
$ws.Range("I23").Value = 41.31896883144262
$ws.Range("J23").Value = 11.3273110035165
$ws.Range("K23").Value = 14.90919627548003
$ws.Range("L23").Value = 10.85659215347193
$ws.Range("M23").Value = 16.45638995256838
$ws.Range("B24").Value = 16.16656893802262
$ws.Range("C24").Value = 6.15153879548502
$ws.Range("E24").Value = 9.932763669789308
$ws.Range("F24").Value = 54.27119945296901
$ws.Range("G24").Value = 3.777208905719282
$ws.Range("I24").Value = 41.03334429560653
$ws.Range("J24").Value = 11.32181182700429
$ws.Range("K24").Value = 14.79717941273187
$ws.Range("L24").Value = 10.89307286858099
$ws.Range("M24").Value = 16.48346595592707
$ws.Range("B25").Value = 16.04303118248624
$ws.Range("C25").Value = 5.810160408386474
$ws.Range("E25").Value = 9.979999837921753
$ws.Range("F25").Value = 53.7973239460194
$ws.Range("G25").Value = 3.782495113820483
$ws.Range("I25").Value = 40.73732497508423
$ws.Range("J25").Value = 11.31577150574073
$ws.Range("K25").Value = 14.72425186525396
$ws.Range("L25").Value = 10.93875911531207
$ws.Range("M25").Value = 16.52474721210239
